$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Row 2 (Treatment / Campaign / MDA, age 5-15): the coverage value of 0.6 is kept
# only for 2018, 2020, 2022, 2024 (H2, J2, L2, N2); all later years are cleared.
$ws.Range("P2:AD2").ClearContents()

# Insert a brand-new row at row 3, shifting the existing rows 3-7 down to 4-8.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 (Treatment / Campaign / MDA, age 2-15) with its
# coverage of 0.8 for 2026, 2028, 2030, 2032, 2034, 2036, 2038, 2040.
$ws.Range("A3").Value = "All"
$ws.Range("B3").Value = "Treatment"
$ws.Range("C3").Value = "Campaign"
$ws.Range("D3").Value = "MDA"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 15
$ws.Range("P3").Value = 0.8
$ws.Range("R3").Value = 0.8
$ws.Range("T3").Value = 0.8
$ws.Range("V3").Value = 0.8
$ws.Range("X3").Value = 0.8
$ws.Range("Z3").Value = 0.8
$ws.Range("AB3").Value = 0.8
$ws.Range("AD3").Value = 0.8

# Row 4 (was row 3 before the insert, Treatment / Campaign / MDA, age 15-50) gains a
# coverage value of 0.5 for the same later years.
$ws.Range("P4").Value = 0.5
$ws.Range("R4").Value = 0.5
$ws.Range("T4").Value = 0.5
$ws.Range("V4").Value = 0.5
$ws.Range("X4").Value = 0.5
$ws.Range("Z4").Value = 0.5
$ws.Range("AB4").Value = 0.5
$ws.Range("AD4").Value = 0.5

# Row 5 (was row 4 before the insert, Treatment / Campaign / MDA, age 50-65) also gains
# a coverage value of 0.5 for the same later years.
$ws.Range("P5").Value = 0.5
$ws.Range("R5").Value = 0.5
$ws.Range("T5").Value = 0.5
$ws.Range("V5").Value = 0.5
$ws.Range("X5").Value = 0.5
$ws.Range("Z5").Value = 0.5
$ws.Range("AB5").Value = 0.5
$ws.Range("AD5").Value = 0.5

# Match the author's final view state: zoomed to 120% with AD2 selected.
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 120
$ws.Range("AD2").Select() | Out-Null
